$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with new sensor readings (custom accuracy refresh)
$ws.Range("A2").Value = 45178.50694444445
$ws.Range("B2").Value = 9.686
$ws.Range("C2").Value = 6.9
$ws.Range("D2").Value = 2.907
$ws.Range("E2").Value = 21.429
$ws.Range("F2").Value = 15.753
$ws.Range("G2").Value = 7.426
$ws.Range("H2").Value = 21.228
$ws.Range("I2").Value = 12.148
$ws.Range("J2").Value = 5.042
$ws.Range("K2").Value = 6.422
$ws.Range("L2").Value = 8.574
$ws.Range("M2").Value = 9.294
$ws.Range("N2").Value = 2.529
$ws.Range("O2").Value = 7.895
$ws.Range("P2").Value = 10.528
$ws.Range("Q2").Value = 7.369
$ws.Range("R2").Value = 2.31
$ws.Range("S2").Value = 0.826
$ws.Range("T2").Value = 112.771
$ws.Range("U2").Value = 21.717
$ws.Range("V2").Value = 7.288
$ws.Range("W2").Value = 13.637
$ws.Range("X2").Value = 7.423
$ws.Range("Y2").Value = 1.536
$ws.Range("Z2").Value = 12.849
$ws.Range("AA2").Value = 6.437
$ws.Range("AB2").Value = 5.959
$ws.Range("AC2").Value = 6.92
$ws.Range("AD2").Value = 9.375999999999999
$ws.Range("AE2").Value = 2.459
$ws.Range("AF2").Value = 19.236
$ws.Range("AG2").Value = 3.642
$ws.Range("AH2").Value = 9.109999999999999

$ws.Range("A3").Value = 45178.51388888889
$ws.Range("B3").Value = 3.593
$ws.Range("C3").Value = 2.541
$ws.Range("D3").Value = 1.234
$ws.Range("E3").Value = 8.259
$ws.Range("F3").Value = 5.758
$ws.Range("G3").Value = 2.67
$ws.Range("H3").Value = 13.173
$ws.Range("I3").Value = 4.592
$ws.Range("J3").Value = 2.006
$ws.Range("K3").Value = 2.139
$ws.Range("L3").Value = 3.327
$ws.Range("M3").Value = 3.642
$ws.Range("N3").Value = 0.974
$ws.Range("O3").Value = 3.008
$ws.Range("P3").Value = 3.992
$ws.Range("Q3").Value = 3.04
$ws.Range("R3").Value = 1.1
$ws.Range("S3").Value = 0.342
$ws.Range("T3").Value = 38.488
$ws.Range("U3").Value = 8.605
$ws.Range("V3").Value = 2.776
$ws.Range("W3").Value = 5.296
$ws.Range("X3").Value = 2.798
$ws.Range("Y3").Value = 0.499
$ws.Range("Z3").Value = 7.362
$ws.Range("AA3").Value = 2.452
$ws.Range("AB3").Value = 2.371
$ws.Range("AC3").Value = 2.769
$ws.Range("AD3").Value = 3.577
$ws.Range("AE3").Value = 1.098
$ws.Range("AF3").Value = 12.673
$ws.Range("AG3").Value = 1.293
$ws.Range("AH3").Value = 3.473

$ws.Range("A4").Value = 45178.52083333334
$ws.Range("B4").Value = 16.633
$ws.Range("C4").Value = 12.403
$ws.Range("D4").Value = 1.271
$ws.Range("E4").Value = 36.537
$ws.Range("F4").Value = 29.489
$ws.Range("G4").Value = 12.958
$ws.Range("H4").Value = 46.834
$ws.Range("I4").Value = 20.304
$ws.Range("J4").Value = 9.093999999999999
$ws.Range("K4").Value = 12.975
$ws.Range("L4").Value = 14.674
$ws.Range("M4").Value = 15.628
$ws.Range("N4").Value = 4.225
$ws.Range("O4").Value = 13.159
$ws.Range("P4").Value = 18.536
$ws.Range("Q4").Value = 11.291
$ws.Range("R4").Value = 0.867
$ws.Range("S4").Value = 0.629
$ws.Range("T4").Value = 192.947
$ws.Range("U4").Value = 36.673
$ws.Range("V4").Value = 12.146
$ws.Range("W4").Value = 24.445
$ws.Range("X4").Value = 12.905
$ws.Range("Y4").Value = 1.769
$ws.Range("Z4").Value = 23.995
$ws.Range("AA4").Value = 10.728
$ws.Range("AB4").Value = 9.574
$ws.Range("AC4").Value = 11.256
$ws.Range("AD4").Value = 15.46
$ws.Range("AE4").Value = 0.707
$ws.Range("AF4").Value = 42.554
$ws.Range("AG4").Value = 6.726
$ws.Range("AH4").Value = 15.185

$ws.Range("A5").Value = 45178.52777777778
$ws.Range("B5").Value = 19.55
$ws.Range("C5").Value = 14.63
$ws.Range("D5").Value = 1.18
$ws.Range("E5").Value = 42.83
$ws.Range("F5").Value = 34.87
$ws.Range("G5").Value = 15.28
$ws.Range("H5").Value = 59.46
$ws.Range("I5").Value = 23.8
$ws.Range("J5").Value = 10.72
$ws.Range("K5").Value = 15.47
$ws.Range("L5").Value = 17.19
$ws.Range("M5").Value = 18.28
$ws.Range("N5").Value = 4.95
$ws.Range("O5").Value = 15.41
$ws.Range("P5").Value = 21.85
$ws.Range("Q5").Value = 13.07
$ws.Range("R5").Value = 0.72
$ws.Range("S5").Value = 0.67
$ws.Range("T5").Value = 227.31
$ws.Range("U5").Value = 43.11
$ws.Range("V5").Value = 14.23
$ws.Range("W5").Value = 28.9
$ws.Range("X5").Value = 15.19
$ws.Range("Y5").Value = 2.05
$ws.Range("Z5").Value = 29.56
$ws.Range("AA5").Value = 12.57
$ws.Range("AB5").Value = 11.15
$ws.Range("AC5").Value = 13.12
$ws.Range("AD5").Value = 18.1
$ws.Range("AE5").Value = 0.52
$ws.Range("AF5").Value = 54.21
$ws.Range("AG5").Value = 7.96
$ws.Range("AH5").Value = 17.79

# Remove the last data row (row 6) - data trimmed to 1000 records overall
$ws.Rows.Item(6).Delete()

# Adjust a handful of column widths slightly
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 6.166666666666667
